$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "DONE" formatting (style) from C6 and apply it to D6:D9,
# then set their values to "DONE" to mark the report column as finished.
$ws.Range("C6").Copy()
$ws.Range("D6:D9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D6:D9").Value = "DONE"

# Update the active selection to reflect where the edit was made.
[void]$ws.Range("C6:D6").Select()
